$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A60").Value = 0.1
$ws.Range("A61").Value = 0.2
$ws.Range("A62").Value = 0.3
$ws.Range("A63").Value = 1.23456789012345
$ws.Range("A64").Value = -66.666666666666
$ws.Range("A65").Value = 12.1
$ws.Range("A66").Value = 100.1
